# Update the weekly Fruta/Hortaliza (Pomelo) price records.
# Rows 3-6 hold the last few weeks of data; this edit rotates the data so the
# newest week moves up and the oldest week moves down, as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 new values
$ws.Range("D3").Value = 44229
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11364
$ws.Range("S3").Value = 812

# Row 4 new values
$ws.Range("D4").Value = 44210
$ws.Range("M4").Value = 70
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 11000
$ws.Range("P4").Value = 10357
$ws.Range("S4").Value = 740

# Row 5 new values
$ws.Range("D5").Value = 44216
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 11545
$ws.Range("S5").Value = 825

# Row 6 new values
$ws.Range("D6").Value = 44172
$ws.Range("M6").Value = 90
$ws.Range("N6").Value = 8500
$ws.Range("O6").Value = 9000
$ws.Range("P6").Value = 8806
$ws.Range("S6").Value = 629
